$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column F ("day") filled down first for the existing rows ---
$ws.Range("F1").Value = "day"
$ws.Range("F1").WrapText = $true
$ws.Range("F2").Value = "Monday"
$ws.Range("F3").Value = "Tuesday"

# --- Column E ("Gender") filled down next for the existing rows ---
$ws.Range("E1").Value = "Gender"
$ws.Range("E1").WrapText = $true
$ws.Range("E2").Value = "Female"
$ws.Range("E3").Value = "Male"

# --- New row 4 (Harish's record), entered out of A-to-F order ---
$ws.Range("B4").Value = "kinnu@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:kinnu@gmail.com")
$ws.Range("C4").Value = 1234567893
$ws.Range("D4").Value = "tirupati"
$ws.Range("F4").Value = "Wednesday"
$ws.Range("E4").Value = "Male"

# --- New row 5 ---
$ws.Range("A5").Value = "k"
$ws.Range("A5").WrapText = $true
$ws.Range("B5").Value = "l@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:l@gmail.com")
$ws.Range("D5").Value = "tirupatiswamy"
$ws.Range("F5").Value = "Thursday"
$ws.Range("E5").Value = "Female"
$ws.Range("C5").Value = 1234567894
$ws.Range("C5").WrapText = $true

# --- A4 typed in last ---
$ws.Range("A4").Value = "Harish"
$ws.Range("A4").WrapText = $true

# --- Selection moved to A4, matching the saved cursor position ---
$ws.Range("A4").Select()
